$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Swap the presentation theme's colour scheme from "Integral" to the
#    default "Office Theme" colours (dk1, lt1, dk2, lt2, accent1-6, hlink,
#    folHlink). RGB values are expressed as PowerPoint's BGR-packed integer
#    (B*65536 + G*256 + R), matching the COM RGBColor.RGB convention.
# ---------------------------------------------------------------------------
$m = $p.SlideMaster
$cs = $m.Theme.ThemeColorScheme

$cs.Item(1).RGB  = 0          # dk1      000000
$cs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$cs.Item(3).RGB  = 6968388    # dk2      44546A
$cs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$cs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$cs.Item(6).RGB  = 3243501    # accent2  ED7D31
$cs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$cs.Item(8).RGB  = 49407      # accent4  FFC000
$cs.Item(9).RGB  = 12874308   # accent5  4472C4
$cs.Item(10).RGB = 4697456    # accent6  70AD47
$cs.Item(11).RGB = 12673797   # hlink    0563C1
$cs.Item(12).RGB = 7491477    # folHlink 954F72

# ---------------------------------------------------------------------------
# 2) Re-style the table on slide 6 to use the built-in table style
#    {4EFF7CA6-0099-4C72-8C3B-CADCC976E81B} instead of the custom
#    {FE431121-B6B8-4EFB-AFC1-615D2BA92AA7} style.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(6)
$shp = $s.Shapes.Item(2)
$shp.Table.ApplyStyle("{4EFF7CA6-0099-4C72-8C3B-CADCC976E81B}")
